# Rearrange row order for template instruction rows in all three sheets.
# Pattern: (names-row, numbers-row, types-row) -> (numbers-row, types-row, names-row)

$wb = $excel.ActiveWorkbook

# ---- Sheet1 ----
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A5").Value = 1
$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 4

$ws1.Range("A6").Value = "iri"
$ws1.Range("B6").Value = "xsd:int"
$ws1.Range("C6").Value = "iri"
$ws1.Range("D6").Value = "iri"

$ws1.Range("A7").Value = "class"
$ws1.Range("B7").Value = "cardinality"
$ws1.Range("C7").Value = "property"
$ws1.Range("D7").Value = "range"

$ws1.Range("B6").Style = $ws1.Range("B7").Style

# ---- Sheet1_2 ----
$ws2 = $wb.Worksheets.Item("Sheet1_2")

$ws2.Range("A8").Value = 4
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = 3
$ws2.Range("D8").Value = 0
$ws2.Range("E8").Value = 1

$ws2.Range("A9").Value = "iri"
$ws2.Range("B9").Value = "xsd:int"
$ws2.Range("C9").Value = "iri"
$ws2.Range("D9").ClearContents()
$ws2.Range("E9").Value = "iri"

$ws2.Range("A10").Value = "class"
$ws2.Range("B10").Value = "cardinality"
$ws2.Range("C10").Value = "property"
$ws2.Range("E10").Value = "range"

$ws2.Range("B9").Style = $ws2.Range("B10").Style

# ---- Sheet3 ----
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("A5").Value = 1
$ws3.Range("A6").Value = "iri+"
$ws3.Range("A7").Value = "classes"
$ws3.Range("B6").Value = ""
